# Apply cryptos list update (prices + 1h volume %) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.955.13"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "2.363.46"
$ws.Range("E3").Value = "  +1.81%  "
$ws.Range("E4").Value = "  +0.00%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "302.41"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.30%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "95.55"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  -0.45%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.488"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.82%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "34.05"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.51%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.124"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +3.29%  "
$ws.Range("E12").Value = "  +0.10%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "18.41"
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.71"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.15%  "
$ws.Range("D15").Value = "2.731.68"
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("D16").Value = "2.366.95"
$ws.Range("E16").Value = "  +1.93%  "
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").Value = "42.921.02"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("E19").Value = "  -2.18%  "
$ws.Range("E20").Value = "  +1.54%  "
$ws.Range("D21").Value = "0.0₃0885"
$ws.Range("E21").Value = "  -0.53%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "68.03"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.34%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "235.26"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("E24").Value = "  -4.65%  "
$ws.Range("E25").Value = "  +1.13%  "
$ws.Range("E26").Value = "  -0.02%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "24.50"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.30%  "
$ws.Range("E28").Value = "  +0.58%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "9.29"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +1.88%  "
$ws.Range("E30").Value = "  -0.54%  "
$ws.Range("E31").Value = "  -0.06%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "5.00"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.08%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "17.45"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.34%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "130.93"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -10.81%  "
$ws.Range("E35").Value = "  +2.14%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.104"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +3.73%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.84"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +2.16%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "4.33"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.85%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.82"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +3.12%  "
$ws.Range("E40").Value = "  -2.06%  "
$ws.Range("E41").Value = "  -0.60%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "21.02"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -4.07%  "
$ws.Range("D43").Value = "1.931.70"
$ws.Range("E43").Value = "  +0.44%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.0279"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +2.86%  "
$ws.Range("E46").Value = "  -9.04%  "
$ws.Range("E47").Value = "  -1.33%  "
$ws.Range("D48").Value = "2.589.73"
$ws.Range("E48").Value = "  +1.44%  "
$ws.Range("E49").Value = "  +2.20%  "
$ws.Range("E50").Value = "  +1.88%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "71.48"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -1.09%  "
